# Refresh "想去人数" (column F) and two "Cover" image URLs (column I)
# on every worksheet, matching the upstream data pull baked into
# gh-pages output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 12571
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 1295
$ws.Range("F13").Value = 5484
$ws.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202407/NSQarDy41720678771123.jpeg"
$ws.Range("F14").Value = 912
$ws.Range("F15").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 1435
$ws.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202407/CtmqI8ub1720675857290.png"
$ws.Range("F19").Value = 0
$ws.Range("F21").Value = 1038
$ws.Range("F23").Value = 882
$ws.Range("F26").Value = 734
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 2063
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 47
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 4439
$ws.Range("F38").Value = 0
$ws.Range("F39").Value = 133
$ws.Range("F41").Value = 646
$ws.Range("F42").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("F49").Value = 193

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F6").Value = 2
$ws.Range("F7").Value = 39
$ws.Range("F11").Value = 77
$ws.Range("F17").Value = 6
$ws.Range("F20").Value = 6
$ws.Range("F21").Value = 13
$ws.Range("F22").Value = 76
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 3
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("F30").Value = 1

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 6539
$ws.Range("F3").Value = 0

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 601
$ws.Range("F3").Value = 0
$ws.Range("F5").Value = 7006
$ws.Range("F6").Value = 141
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 12951
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 1295
$ws.Range("F13").Value = 0
$ws.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202407/NSQarDy41720678771123.jpeg"
$ws.Range("F14").Value = 912
$ws.Range("F15").Value = 365
$ws.Range("F16").Value = 194
$ws.Range("F17").Value = 1435
$ws.Range("I17").Value = "//i2.hdslb.com/bfs/openplatform/202407/CtmqI8ub1720675857290.png"
$ws.Range("F18").Value = 359
$ws.Range("F20").Value = 1038
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 3018
$ws.Range("F25").Value = 5
$ws.Range("F26").Value = 258
$ws.Range("F27").Value = 2063
$ws.Range("F28").Value = 109
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 1687
$ws.Range("F33").Value = 149
$ws.Range("F34").Value = 47
$ws.Range("F35").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 193

